$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(-0.6229298114776611, -4.338823318481445, 0.1601994037628173)
    ,@(-0.3721696436405182, -3.498425245285034, 0.6684392690658569)
    ,@(0.1207986027002334, -1.66308057308197, 0.5348125100135803)
    ,@(-0.5253441333770752, -1.381013631820679, 0.5633704662322998)
    ,@(0.5025894045829773, 2.330298185348511, 0.1441642045974731)
    ,@(0.8329147696495056, 6.17554235458374, -1.043815970420837)
    ,@(0.0348193198442459, 5.04437255859375, 0.868803083896637)
    ,@(-0.2368630021810531, 1.756695747375489, 0.40943244099617)
    ,@(-0.0340557359158992, 2.139860868453979, 0.6606507897377014)
    ,@(0.1806633919477462, 1.172861218452454, 0.4958698749542236)
    ,@(0.0870482996106147, -3.714518785476685, -0.6794348955154419)
    ,@(0.1391245573759079, -6.393468379974365, -4.738176345825195)
    ,@(1.760971784591675, -1.517847418785095, 1.018465042114258)
    ,@(-0.0085521135479211, -4.238488674163818, -0.0195476878434419)
    ,@(0.4407392740249634, -3.014314651489258, -0.6562219858169556)
    ,@(-0.1287398487329483, 2.51019811630249, -0.0731511116027832)
    ,@(0.2913827300071716, 6.333756446838379, 0.9952521920204164)
    ,@(-1.667662143707275, 4.299574851989746, 0.0829249545931816)
    ,@(1.54808521270752, 2.839301586151123, 1.368032693862915)
    ,@(0.0797179117798805, 3.36098051071167, 0.3753767013549804)
    ,@(0.4167627990245819, 1.646892666816711, 0.836885392665863)
    ,@(0.5007568001747131, -1.671174645423889, -0.2092213481664657)
    ,@(0.0099265603348612, -3.577379703521729, -0.4952589869499206)
    ,@(-1.279915452003479, -6.508005619049072, -0.1085812970995903)
    ,@(-0.9758572578430176, -1.915215253829956, 1.272432327270508)
    ,@(0.0806342139840126, -1.140027284622192, 0.0485637858510017)
    ,@(0.8413141369819641, -0.3859141170978546, 0.5590944290161133)
    ,@(-0.1750128865242004, 2.383749008178711, -0.1401935666799545)
    ,@(0.4100432991981506, 2.355190992355347, 0.1059851199388504)
    ,@(-0.6436992287635803, 0.6840163469314575, -0.0784961804747581)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
